## Undo Jason's overwrite of class materials starter code/slides.
##
## 1) Slide 1 "TextBox 2": the attendance-password line changes from the
##    literal word "bigo" to a blanked-out "__________" placeholder.
## 2) The "Date Placeholder" cached field text on the slide master and on
##    every slide layout gets rolled back from "Sunday, September 3, 2023"
##    to "Tuesday, February 7, 2023" (9/3/2023 -> 2/7/2023).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1, shape "TextBox 2": replace the "bigo" paragraph text.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 2" -and $shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if ($para.Text -eq "bigo") {
                $para.Text = "__________"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide master date placeholder.
# ---------------------------------------------------------------------
$masterShapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $masterShapes.Count; $i++) {
    $shp = $masterShapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Sunday, September 3, 2023") {
        $shp.TextFrame.TextRange.Text = "Tuesday, February 7, 2023"
    }
}

# ---------------------------------------------------------------------
# 3) Each slide layout's date placeholder.
# ---------------------------------------------------------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Sunday, September 3, 2023") {
            $shp.TextFrame.TextRange.Text = "Tuesday, February 7, 2023"
        }
    }
}
